# "Generate Report for Handoff"
# The workbook has 3 sheets: Overview, zh-cn, de-de.
# Each sheet currently tracks two handoff entries (rows 2 and 3); the
# second entry (row 3, the "d26906ea..." file) is removed and the
# remaining entry's status/date fields are refreshed to reflect a new
# handoff ("Ready for handoff" / updated timestamps).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-36-18 10:36:27"

$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b06cbdde86c553aa18aac93e4a05aba7b4e82a2/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-18 10:36:25"

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b06cbdde86c553aa18aac93e4a05aba7b4e82a2/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b06cbdde86c553aa18aac93e4a05aba7b4e82a2/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79cf11268a41327a24bc96bf09c0240c63d80d2b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.zh-cn.xlf", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/dfae951f456dcfefa41b571bcf0d2d70bdcc068c/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d4f8436d7d0e07cb5933f7a4551a911d1a5fbac2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.zh-cn.xlf", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.zh-cn.xlf")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-18 10:36:27"

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b06cbdde86c553aa18aac93e4a05aba7b4e82a2/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b06cbdde86c553aa18aac93e4a05aba7b4e82a2/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/de0ce16be7fe4a91ac35aec36e48cab149085621/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.de-de.xlf", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ee774bb5ae54785bdf2f7ba951e6ffcc24b921ec/e2e/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/23cb18fadc00617ad646d4886b3dd1eddf2eabc0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.de-de.xlf", "", "", "5ea04b7e-4820-4e83-bf4a-bb73317ee9b6.e68aec360c515e45bc5a235793ba00fb1496747b.de-de.xlf")

Write-Host "Report generated for handoff"
